# Updates for v1.2. Introduction of inference tool guidelines.
#
# 1. Insert a new "Repertoire" worksheet between "Submission" and "Inferences".
# 2. Submission sheet: insert a blank row (no new content) so the lower block
#    shifts down by one row.
# 3. Inferences sheet: insert 4 rows before the old row 7 and populate the
#    new guidance/subject/genotype-id cells.
# 4. Genotype sheet: insert 3 rows before the old row 7 and populate the new
#    guidance/genotype-id/subject-id cells.
# 5. Tool Settings sheet is untouched.

$wb = $excel.ActiveWorkbook

# --- 1. Add the "Repertoire" worksheet, positioned before "Inferences" ---
$inferencesSheet = $wb.Worksheets.Item("Inferences")
$repertoire = $wb.Worksheets.Add($inferencesSheet)
$repertoire.Name = "Repertoire"

$repertoire.Range("B3").Value = "Repertoire"
$repertoire.Range("B3").Font.Bold = $true
$repertoire.Range("B4").Font.Bold = $true
$repertoire.Range("B5").Value = "Please provide details of the repertoire from which inferences were made, and the methods by which it was constructed."
$repertoire.Range("B8").Value = "{{Repertoire:properties!completed_by,type+Response}}"

$repertoire.Range("G8").Select()

# --- 2. Submission sheet: shift rows 7/10/12 down by one ---
$submission = $wb.Worksheets.Item("Submission")
$submission.Range("A7:A10").EntireRow.Insert()
$submission.Range("B11:G13").Select()

# --- 3. Inferences sheet: shift rows 7+ down by four, add new guidance rows ---
$inferences = $wb.Worksheets.Item("Inferences")
$inferences.Range("A7:A10").EntireRow.Insert()
$inferences.Range("B6").Value = "The table should be used to describe sequences inferred from a single invididual, in a single genotype. If inferences are made from multiple individuals and/or multiple genotypes,"
$inferences.Range("B7").Value = "please create a separate table on this tab for each genotyope from which inferences are made, and cross-reference to separate genotype tables on the genotype tab."
$inferences.Range("B9").Value = "Subject id:"
$inferences.Range("B10").Value = "Genotype id (as listed on the Genotype tab):"
$inferences.Range("B9").Select()

# --- 4. Genotype sheet: shift rows 7+ down by three, add new guidance rows ---
$genotype = $wb.Worksheets.Item("Genotype")
$genotype.Range("A7:A9").EntireRow.Insert()
$genotype.Range("B6").Value = "If inferences are made from multiple genotypes, please create separate tables for each genotype, and cross-reference to the inferences made from each one on the Inferences tab."
$genotype.Range("B8").Value = "Genotype Id:"
$genotype.Range("B9").Value = "Subject Id (as listed on Inferences tab):"
$genotype.Range("B3:B6").Select()

# --- 5. Tool Settings sheet: untouched ---

# Leave the "Repertoire" tab active/selected, matching the committed state.
$repertoire.Activate()
